# Update the cryptocurrency price (D) and volume-change (E) columns
# with the latest scraped values, matching the GitHub Actions refresh job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Force the cell to keep storing a literal string (e.g. '0.999' or
    # '560.54') instead of letting Excel auto-convert it to a float,
    # then restore the cell's original style so no formatting changes.
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '64.688.47'
$ws.Range("E2").Value = '  +2.99%  '
Set-TextValue $ws.Range("D3") '3.086.91'
$ws.Range("E3").Value = '  +1.53%  '
Set-TextValue $ws.Range("D4") '0.999'
$ws.Range("E4").Value = '  -0.24%  '
Set-TextValue $ws.Range("D5") '560.54'
$ws.Range("E5").Value = '  +2.22%  '
Set-TextValue $ws.Range("D6") '146.64'
$ws.Range("E6").Value = '  +6.95%  '
Set-TextValue $ws.Range("D7") '0.999'
$ws.Range("E7").Value = '  -0.15%  '
Set-TextValue $ws.Range("D8") '3.087.31'
$ws.Range("E8").Value = '  +1.79%  '
Set-TextValue $ws.Range("D9") '0.504'
$ws.Range("E9").Value = '  +1.39%  '
Set-TextValue $ws.Range("D10") '0.154'
$ws.Range("E10").Value = '  +4.00%  '
Set-TextValue $ws.Range("D11") '6.14'
$ws.Range("E11").Value = '  -1.03%  '
Set-TextValue $ws.Range("D12") '0.473'
$ws.Range("E12").Value = '  +5.93%  '
Set-TextValue $ws.Range("D13") '0.0000229'
$ws.Range("E13").Value = '  +2.27%  '
Set-TextValue $ws.Range("D14") '35.40'
$ws.Range("E14").Value = '  +2.86%  '
Set-TextValue $ws.Range("D15") '3.572.57'
$ws.Range("E15").Value = '  +0.78%  '
Set-TextValue $ws.Range("D16") '64.613.45'
$ws.Range("E16").Value = '  +2.66%  '
Set-TextValue $ws.Range("D17") '3.079.47'
$ws.Range("E17").Value = '  +0.93%  '
Set-TextValue $ws.Range("D18") '0.110'
$ws.Range("E18").Value = '  +1.77%  '
Set-TextValue $ws.Range("D19") '6.82'
$ws.Range("E19").Value = '  +2.29%  '
Set-TextValue $ws.Range("D20") '482.20'
$ws.Range("E20").Value = '  +1.08%  '
Set-TextValue $ws.Range("D21") '14.01'
$ws.Range("E21").Value = '  +3.92%  '
Set-TextValue $ws.Range("D22") '0.681'
$ws.Range("E22").Value = '  +1.92%  '
Set-TextValue $ws.Range("D23") '7.60'
$ws.Range("E23").Value = '  +6.90%  '
Set-TextValue $ws.Range("D24") '13.81'
$ws.Range("E24").Value = '  +12.01%  '
Set-TextValue $ws.Range("D25") '81.99'
$ws.Range("E25").Value = '  +1.83%  '
Set-TextValue $ws.Range("D26") '1.00'
$ws.Range("E26").Value = '  +0.17%  '
Set-TextValue $ws.Range("D27") '2.82'
$ws.Range("E27").Value = '  +2.93%  '
Set-TextValue $ws.Range("D28") '8.20'
$ws.Range("E28").Value = '  +5.28%  '
Set-TextValue $ws.Range("D29") '2.09'
$ws.Range("E29").Value = '  +6.85%  '
Set-TextValue $ws.Range("D30") '0.998'
$ws.Range("E30").Value = '  -0.41%  '
Set-TextValue $ws.Range("D31") '26.32'
$ws.Range("E31").Value = '  +2.23%  '
Set-TextValue $ws.Range("D32") '1.15'
$ws.Range("E32").Value = '  +1.15%  '
Set-TextValue $ws.Range("D33") '2.51'
$ws.Range("E33").Value = '  +5.52%  '
Set-TextValue $ws.Range("D34") '5.62'
$ws.Range("E34").Value = '  +0.15%  '
Set-TextValue $ws.Range("D35") '6.24'
$ws.Range("E35").Value = '  +5.67%  '
Set-TextValue $ws.Range("D36") '54.92'
$ws.Range("E36").Value = '  -0.13%  '
Set-TextValue $ws.Range("D37") '463.87'
$ws.Range("E37").Value = '  +1.27%  '
Set-TextValue $ws.Range("D38") '3.03'
$ws.Range("E38").Value = '  +20.55%  '
Set-TextValue $ws.Range("D39") '0.0834'
$ws.Range("E39").Value = '  +3.51%  '
Set-TextValue $ws.Range("D40") '0.0407'
$ws.Range("E40").Value = '  +4.42%  '
Set-TextValue $ws.Range("D41") '2.974.69'
$ws.Range("E41").Value = '  -4.06%  '
Set-TextValue $ws.Range("D42") '8.30'
$ws.Range("E42").Value = '  +1.62%  '
Set-TextValue $ws.Range("D43") '0.115'
$ws.Range("E43").Value = '  -0.24%  '
Set-TextValue $ws.Range("D44") '27.94'
$ws.Range("E44").Value = '  +1.09%  '
Set-TextValue $ws.Range("D45") '0.265'
$ws.Range("E45").Value = '  +6.71%  '
Set-TextValue $ws.Range("D46") '2.17'
$ws.Range("E46").Value = '  +7.81%  '
Set-TextValue $ws.Range("D47") '1.00'
$ws.Range("E47").Value = '  +0.07%  '
Set-TextValue $ws.Range("D48") '0.113'
$ws.Range("E48").Value = '  +3.76%  '
Set-TextValue $ws.Range("D49") '120.65'
$ws.Range("E49").Value = '  +4.02%  '
Set-TextValue $ws.Range("D50") '0.0₃0520'
$ws.Range("E50").Value = '  +4.06%  '
Set-TextValue $ws.Range("D51") '2.10'
$ws.Range("E51").Value = '  +3.02%  '
